$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column D: header "reverses" (bold, centered/top, thin left+right border)
# and data row "no" underneath (plain, like the other data cells).
$ws.Range("D1").Value = "reverses"
$ws.Range("D2").Value = "no"

# Build the header's border style (thin left + thin right only) on a scratch
# cell first so it can be copied onto D1 as a single formatting operation -
# this keeps the font/alignment/border combination together in one new
# cell style instead of being fragmented across several intermediate ones.
$scratch = $ws.Range("F1")
$scratch.Value = "scratch"
$scratch.Font.Bold = $true
$scratch.HorizontalAlignment = -4108   # xlCenter
$scratch.VerticalAlignment = -4160     # xlTop
$scratch.Borders.Item(7).LineStyle = 1   # xlEdgeLeft   = xlContinuous
$scratch.Borders.Item(10).LineStyle = 1  # xlEdgeRight  = xlContinuous

$scratch.Copy()
$ws.Range("D1").PasteSpecial(-4122)   # xlPasteFormats

$scratch.Clear()

# Leave the selection on D2, matching the saved workbook state.
$ws.Range("D2").Select() | Out-Null
